# Finish RPA script updates
# Populate the oct_2020 (column H) figures that were missing for SFY 2021,
# and roll them into the SFY 2021 Total (column Q) for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 228
    4  = 42
    5  = 270
    6  = 1
    7  = 1
    8  = 82
    9  = 8
    10 = 10
    13 = 30
    14 = 132
    15 = 50
    17 = 14
    18 = 17
    20 = 16
    21 = 31
    22 = 3
    23 = 2
    24 = 4
    25 = 23
    26 = 46
    27 = 5
    28 = 2
    29 = 33
    31 = 3
    33 = 2
    34 = 11
    36 = 28
    40 = 3
    41 = 1
    43 = 1
    46 = 3
    47 = 12
}

# Row 15 is the only row where the SFY 2021 Total (Q) is NOT refreshed.
$skipTotalRows = @(15)

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("H$row").Value = $value
    if ($skipTotalRows -notcontains $row) {
        $ws.Range("Q$row").Value = $value
    }
}
